# Updated symbol list on Tue Dec 13 12:00:16 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price values (column D) for rows where the coin has a real price
$priceUpdates = @(
    @{Row=2;  Value="268.44"},
    @{Row=3;  Value="21.41"},
    @{Row=4;  Value="6.245"},
    @{Row=5;  Value="0.06162"},
    @{Row=6;  Value="3.567"},
    @{Row=7;  Value="6.557"},
    @{Row=8;  Value="1.372"},
    @{Row=9;  Value="0.8227"},
    @{Row=10; Value="0.01350"},
    @{Row=11; Value="0.1549"},
    @{Row=12; Value="0.08175"},
    @{Row=13; Value="0.03297"},
    @{Row=14; Value="0.03198"},
    @{Row=15; Value="0.09295"},
    @{Row=16; Value="3.750"},
    @{Row=17; Value="0.001657"},
    @{Row=18; Value="0.04685"},
    @{Row=19; Value="0.006323"},
    @{Row=20; Value="0.005771"},
    @{Row=21; Value="0.001068"},
    @{Row=22; Value="0.0001499"},
    @{Row=24; Value="2.318"},
    @{Row=28; Value="0.0001617"},
    @{Row=40; Value="0.04659"},
    @{Row=41; Value="0.007018"},
    @{Row=42; Value="0.003896"},
    @{Row=44; Value="0.01186"},
    @{Row=45; Value="0.00006031"},
    @{Row=46; Value="0.0009890"},
    @{Row=48; Value="0.7813"},
    @{Row=49; Value="0.002439"},
    @{Row=51; Value="0.01239"}
)

foreach ($update in $priceUpdates) {
    $cell = $ws.Cells.Item($update.Row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $update.Value
}

# Updated hour (column G) from 11 to 12 for every data row 2..51
for ($row = 2; $row -le 51; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $cell.NumberFormat = "@"
    $cell.Value = "12"
}
